# Update the "save to libraries.csv" sheet so row 6 picks up the new
# "Batik" library entry (row 7 on the "libraries" sheet) instead of the
# stale #REF! formulas left over from a deleted row.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A6").Formula = "=libraries!A7"
$ws2.Range("B6").Formula = "=libraries!B7"
$ws2.Range("C6").Formula = "=libraries!C7"
$ws2.Range("D6").Formula = "=libraries!D7"

# Make "save to libraries.csv" the active/selected sheet, with the new
# row just below the table selected (mirrors the author re-selecting the
# refreshed range after fixing the formulas).
$ws2.Activate()
$ws2.Range("A20:D20").Select()
